$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Start clean: wipe existing content/formatting ---
$ws.Cells.Clear()

# --- Header row ---
$ws.Range("A1").Value = "Patient ID"
$ws.Range("B1").Value = "Order template"
$ws.Range("C1").Value = "Start date"
$ws.Range("D1").Value = "Start time"
$ws.Range("E1").Value = "End date"
$ws.Range("F1").Value = "End time"

# --- Row 2 : Cardioversion ---
$ws.Range("A2").Value = 1000000008
$ws.Range("B2").Value = "Cardioversion"
$ws.Range("C2").NumberFormat = "h:mm"
$ws.Range("C2").Value = "20/01/2025"
$ws.Range("D2").NumberFormat = "h:mm"
$ws.Range("D2").Value = 0.3125
$ws.Range("E2").NumberFormat = "h:mm"
$ws.Range("E2").Value = "21/01/2025"
$ws.Range("F2").NumberFormat = "h:mm"
$ws.Range("F2").Value = 0.5

# --- Row 3 : Echokardiographie TTE (quote-prefixed patient id/template) ---
$ws.Range("Z1").Value = "'1000000003"
$ws.Range("A3").Value = 1000000003
$ws.Range("Z1").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

$ws.Range("B3").Value = "'Echokardiographie TTE"

$ws.Range("C3").NumberFormat = "h:mm"
$ws.Range("C3").Value = "19/01/2025"
$ws.Range("D3").NumberFormat = "h:mm"
$ws.Range("D3").Value = 0.35416666666666702
$ws.Range("E3").NumberFormat = "h:mm"
$ws.Range("E3").Value = "20/01/2025"
$ws.Range("F3").NumberFormat = "h:mm"
$ws.Range("F3").Value = 0.375

# --- Row 4 : MRT (DCM03) ---
$ws.Range("A4").Value = 1000000421
$ws.Range("B4").Value = "MRT (DCM03)"
$ws.Range("C4").NumberFormat = "h:mm"
$ws.Range("C4").Value = "17/01/2025"
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("D4").Value = 0.39583333333333298
$ws.Range("E4").NumberFormat = "h:mm"
$ws.Range("E4").Value = "19/01/2025"
$ws.Range("F4").NumberFormat = "h:mm"
$ws.Range("F4").Value = 0.58333333333333304

# --- Row 5 : Ultrasound [US] (DCM01) (quote-prefixed patient id/template) ---
$ws.Range("Z1").Value = "'1000000005"
$ws.Range("A5").Value = 1000000005
$ws.Range("Z1").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

$ws.Range("B5").Value = "'Ultrasound [US] (DCM01)"

$ws.Range("C5").NumberFormat = "h:mm"
$ws.Range("C5").Value = "15/01/2025"
$ws.Range("D5").NumberFormat = "h:mm"
$ws.Range("D5").Value = 0.4375
$ws.Range("E5").NumberFormat = "h:mm"
$ws.Range("E5").Value = "25/01/2025"
$ws.Range("F5").NumberFormat = "h:mm"
$ws.Range("F5").Value = 0.625

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 15.28515625
$ws.Columns("B").ColumnWidth = 27.85546875
$ws.Columns("C").ColumnWidth = 19.28515625
$ws.Columns("D").ColumnWidth = 9.28515625
$ws.Columns("E").ColumnWidth = 10.42578125

# --- Selection ---
$ws.Range("E12").Select()
